# Applies the "variogram01.docx" edit:
#  1. Removes the standalone paragraph that only held a red, noProof
#     manual page break (<w:br w:type="page"/>) right before the
#     "Solution" heading paragraph.
#  2. Strips the now-stale <w:lastRenderedPageBreak/> marker from the
#     "Solution" run (it no longer starts a rendered page once the
#     explicit page break above it is gone).
#  3. Relocates the document's singleton "_GoBack" bookmark from its old
#     position (inside the "b. Yes. ..." answer paragraph, right after
#     the run containing "n") to the start of the "Solution" paragraph.

function Find-Paragraph($doc, $text) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $para = $doc.Paragraphs.Item($i)
        # paragraph text ends with a paragraph mark (13) or a cell/row
        # mark (7); strip those before comparing.
        $t = $para.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $para
        }
    }
    return $null
}

$d = $word.ActiveDocument

# --- Step 1: delete the paragraph that contains only the manual page
#     break immediately preceding the "Solution" heading (it becomes
#     empty once the trailing page-break character is trimmed off).
$solutionPara = Find-Paragraph $d "Solution"
$prevPara = $solutionPara.Previous()
$prevText = $prevPara.Range.Text.TrimEnd([char]13, [char]7, [char]12)
if ($prevText -eq "") {
    $prevPara.Range.Delete()
}

# --- Step 2: re-locate the (possibly re-indexed) "Solution" paragraph
#     and rewrite its range text in place; this preserves the run's
#     character formatting (bold/noProof) while dropping the stale
#     <w:lastRenderedPageBreak/> run child that Find/Delete alone can't
#     touch (it isn't part of the text stream).
$solutionPara = Find-Paragraph $d "Solution"
$solutionPara.Range.Text = "Solution"

# --- Step 3: move the "_GoBack" bookmark to the start of the
#     "Solution" paragraph. Re-adding a bookmark named "_GoBack" moves
#     the existing one (Word treats it as a singleton last-edit marker),
#     so this single call both removes it from its old spot (after the
#     "n" run) and creates it at the new one.
$solutionPara = Find-Paragraph $d "Solution"
$bmRange = $solutionPara.Range.Duplicate
$bmRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRange)
